# =====================================================================
# Adds additional player-scraping data to the workbook:
#   1. New "Player Info" sheet (inserted first)
#   2. "ODI Batting"  sheet: MATCH_CARD_LINK -> MATCH_CODE (bare match code)
#   3. "ODI Bowling"  sheet: MATCH_CARD_LINK -> MATCH_CODE (bare match code)
#   4. New "ODI Batting Extra" sheet (appended last)
# =====================================================================

$wb = $excel.ActiveWorkbook

# Helper: apply the workbook's existing bold/bordered/centered header style
# to a range (this reuses cellXf index 1 from styles.xml instead of
# fabricating a brand-new one).
function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous (thin)
}

# NOTE on worksheet handles in this COM runtime: a worksheet reference
# behaves like a *positional slot* handle rather than a stable object
# identity. Inserting a sheet BEFORE an existing tab shifts every handle
# at/after that slot to point at whatever now occupies it. To stay safe,
# every sheet reference below is re-fetched by name immediately before
# it is used, and always AFTER any `Worksheets.Add(...)` call that could
# have shifted slots.

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet -> inserted as the very first tab
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle $playerInfo.Range("A1:D1")

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3952"
$playerInfo.Range("B2").Value = "Shannon Terry Gabriel"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK column -> MATCH_CODE (bare code)
# ---------------------------------------------------------------------
$battingCodes = @(
    "3905", "3907", "3909", "3939", "3944", "3960", "3961", "3963", "4001",
    "4004", "4017", "4018", "4019", "4040", "4043", "4100", "4101", "4102",
    "4285", "4286", "4291", "4296", "4321", "4325", "4344"
)

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2:D26").NumberFormat = "@"
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $battingSheet.Range("D$row").Value = $battingCodes[$i]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK column -> MATCH_CODE (bare code)
# ---------------------------------------------------------------------
$bowlingCodes = $battingCodes

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2:B26").NumberFormat = "@"
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $bowlingSheet.Range("B$row").Value = $bowlingCodes[$i]
}

# ---------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet -> appended as the very last tab
#    (appending AFTER the last sheet does not shift any existing slots,
#    but we still re-fetch "ODI Bowling" fresh to be safe)
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowlingSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle $extra.Range("A1:F1")

# MATCH_CODE values (column A) are plain text, same codes as above
$extra.Range("A2:A21").NumberFormat = "@"
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $extra.Range("A$row").Value = $battingCodes[$i]
}

# BATTING_POSITION (column B) holds genuine numbers where known
$battingPositions = @{
    2  = 11   # 3960
    4  = 11   # 3963
    5  = 10   # 4001
    7  = 11   # 4017
    9  = 11   # 4019
    10 = 11   # 4040
    11 = 10   # 4043
    12 = 11   # 4100
    13 = 11   # 4101
}
foreach ($row in $battingPositions.Keys) {
    $extra.Range("B$row").Value = $battingPositions[$row]
}

# NUM_4 / NUM_6 (columns C, D) hold text "0" where known
$zeroFourSixRows = @(2, 4, 7, 10, 13)
$extra.Range("C2").NumberFormat = "@"
foreach ($row in $zeroFourSixRows) {
    $extra.Range("C$row").NumberFormat = "@"
    $extra.Range("C$row").Value = "0"
    $extra.Range("D$row").NumberFormat = "@"
    $extra.Range("D$row").Value = "0"
}

# PERCENT_RUNS_OF_TOTAL (column E) — only one known value
$extra.Range("E4").Value = "0.44%"

# MAN_OF_MATCH (column F) — "NO" for the first 13 data rows (rows 2-14),
# left blank for the remaining rows (15-21)
for ($row = 2; $row -le 14; $row++) {
    $extra.Range("F$row").Value = "NO"
}

# ---------------------------------------------------------------------
# Make "Player Info" the active tab (matches activeTab="0" in the target)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Player Info").Activate()
